$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NMA comparison table values with the corrected estimates.
$ws.Range("B2").Value = "0.17 (-4.95;  5.28)"
$ws.Range("F2").Value = "-4.19 (-6.91; -1.48)"

$ws.Range("A3").Value = "-1.43 ( -6.07;  3.20)"
$ws.Range("F3").Value = "-1.16 (-6.28;  3.97)"

$ws.Range("A4").Value = "-1.91 ( -5.59;  1.78)"
$ws.Range("B4").Value = "-0.48 ( -5.74;  4.79)"
$ws.Range("F4").Value = "-2.28 (-4.77;  0.21)"

$ws.Range("A5").Value = "-2.82 ( -8.65;  3.02)"
$ws.Range("B5").Value = "-1.39 ( -8.33;  5.56)"
$ws.Range("C5").Value = "-0.91 ( -6.65;  4.83)"
$ws.Range("F5").Value = "-1.37 (-6.54;  3.79)"

$ws.Range("A6").Value = "-4.48 (-10.32;  1.36)"
$ws.Range("B6").Value = "-3.05 ( -9.99;  3.90)"
$ws.Range("C6").Value = "-2.57 ( -8.31;  3.16)"
$ws.Range("D6").Value = "-1.66 ( -8.97;  5.65)"
$ws.Range("F6").Value = "0.29 (-4.88;  5.46)"

$ws.Range("A7").Value = "-4.19 ( -6.91; -1.48)"
$ws.Range("B7").Value = "-2.76 ( -7.40;  1.88)"
$ws.Range("C7").Value = "-2.28 ( -4.77;  0.21)"
$ws.Range("D7").Value = "-1.37 ( -6.54;  3.79)"
$ws.Range("E7").Value = "0.29 ( -4.88;  5.46)"

# Re-fit column widths so they reflect the new (longer, now-negative) text,
# matching the bestFit recalculation Excel performs automatically.
$ws.Columns.Item(1).ColumnWidth = 16.833333333333332
$ws.Columns.Item(2).ColumnWidth = 16.166666666666668
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 16.166666666666668
$ws.Columns.Item(5).ColumnWidth = 15.5
$ws.Columns.Item(6).ColumnWidth = 16.0

$wb.Save()
